$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.315.62'
$ws.Range("E2").Value = '  +0.02%  '

$ws.Range("D3").Value = '1.876.65'
$ws.Range("E3").Value = '  +0.25%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7108'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.40%  '

$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07996'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.78%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3152'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.39%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.96'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.41%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08267'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.53%  '

$ws.Range("D12").Value = '1.894.25'
$ws.Range("E12").Value = '  +1.53%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.244'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.23%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.48'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.75%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7117'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.09%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.342'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.21%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008534'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.81%  '

$ws.Range("D18").Value = '29.336.48'
$ws.Range("E18").Value = '  +0.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.44'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.70%  '

$ws.Range("D20").Value = '2.144.80'
$ws.Range("E20").Value = '  +1.54%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.49%  '

$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.783'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1553'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.90%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.047'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.48'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.52'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.507'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.414'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.314'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.19%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05366'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.36%  '

$ws.Range("E33").Value = '  -8.74%  '

$ws.Range("E34").Value = '  -0.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7649'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.72%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.182'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.59%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.687'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.51%  '

$ws.Range("E38").Value = '  +0.74%  '

$ws.Range("D39").Value = '1.256.25'
$ws.Range("E39").Value = '  +2.36%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.750'
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.509'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9149'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.30%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '113.01'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.12%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '74.12'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.13%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000133'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.96%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.05%  '

$ws.Range("D47").Value = '2.037.37'
$ws.Range("E47").Value = '  +0.99%  '

$ws.Range("E48").Value = '  +0.41%  '

$ws.Range("E49").Value = '  -0.17%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.455'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.76%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4365'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.14%  '
